$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns keep their original text formatting
# so numeric-looking strings (e.g. "671.02") are not auto-converted to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.691.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.689.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "671.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.26"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.65%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.44%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.17%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000233"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.04"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.720.14"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.682.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.70%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "472.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.647"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "79.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.838.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000127"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.74%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.56%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.69%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.49%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.86"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.50"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.693.82"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.26"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.03%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.02"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.17%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0907"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.82%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.935"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.02"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.29"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.74%  "

# Row 47
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000272"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.37"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.87"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "365.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.63%  "
